$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new batch of 3 price rows (dated 2021-10-07 / serial 44476)
# is inserted right after row 135, pushing the existing rows 136-160 down to
# rows 139-163.
$ws.Rows.Item(136).EntireRow.Insert()
$ws.Rows.Item(136).EntireRow.Insert()
$ws.Rows.Item(136).EntireRow.Insert()

# New row 136: Alcachofa Argentina(o) / Primera
$ws.Range("A136").Value = 2
$ws.Range("B136").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C136").Value = "Coquimbo"
$ws.Range("D136").Value = 44476
$ws.Range("E136").Value = 4
$ws.Range("F136").Value = 100112013
$ws.Range("G136").Value = "Alcachofa"
$ws.Range("H136").Value = "Argentina(o)"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 1100
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 10000
$ws.Range("M136").Value = 9500
$ws.Range("N136").Value = "`$/caja 50 unidades"
$ws.Range("O136").Value = "Provincia de Limarí"
$ws.Range("P136").Value = 190
$ws.Range("Q136").Value = 50
$ws.Range("R136").Value = "Hortaliza"

# New row 137: Alcachofa Española / Primera
$ws.Range("A137").Value = 2
$ws.Range("B137").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C137").Value = "Coquimbo"
$ws.Range("D137").Value = 44476
$ws.Range("E137").Value = 4
$ws.Range("F137").Value = 100112013
$ws.Range("G137").Value = "Alcachofa"
$ws.Range("H137").Value = "Española"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 9000
$ws.Range("L137").Value = 10000
$ws.Range("M137").Value = 9500
$ws.Range("N137").Value = "`$/caja 30 unidades"
$ws.Range("O137").Value = "Provincia de Limarí"
$ws.Range("P137").Value = 317
$ws.Range("Q137").Value = 30
$ws.Range("R137").Value = "Hortaliza"

# New row 138: Alcachofa Madrigal / Primera
$ws.Range("A138").Value = 2
$ws.Range("B138").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C138").Value = "Coquimbo"
$ws.Range("D138").Value = 44476
$ws.Range("E138").Value = 4
$ws.Range("F138").Value = 100112013
$ws.Range("G138").Value = "Alcachofa"
$ws.Range("H138").Value = "Madrigal"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 700
$ws.Range("K138").Value = 7000
$ws.Range("L138").Value = 8000
$ws.Range("M138").Value = 7500
$ws.Range("N138").Value = "`$/caja 40 unidades"
$ws.Range("O138").Value = "Provincia de Limarí"
$ws.Range("P138").Value = 188
$ws.Range("Q138").Value = 40
$ws.Range("R138").Value = "Hortaliza"
